$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a green fill style for the "TAK" column (E) cells ---
# Apply a solid green (RGB 92D050) interior to E3:E7, E10:E14, E17:E19
$ws.Range("E3:E7").Interior.Color = 5296274
$ws.Range("E10:E14").Interior.Color = 5296274
$ws.Range("E17:E19").Interior.Color = 5296274

# --- Fill in the missing "TAK" values for E18 and E19 ---
$ws.Range("E18").Value = "TAK"
$ws.Range("E19").Value = "TAK"

# --- Update the To Do list text ---
# Row 23: "Poprawić obsługę liczb ujemnych..." task is done, replaced by "Konwersja RPN na wynik"
$ws.Range("B23").Value = "Konwersja RPN na wynik"
# Row 24: new task number "2." and a new task description about start/exit screen
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "2."
$ws.Range("B24").Value = "Dodanie ekranu startowego i wyjścia, żeby można było kilkukrotnie wpisać sobie jakieś wyniki"

# --- Update the view: scroll so row 4 is at top, and move selection to A25 ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A25").Select()
